# NMCARS-PART-5232: apply "List 1" style to several (a)/(b)/(c)/(d)/(e)
# numbered paragraphs, add a lastRenderedPageBreak hint before the
# 5232.70 SUBPART heading, tidy the built-in "List 2" style, and add
# the "List 1" / "List 3" / "List 4" list styles (with their linked
# character styles) that the re-macroed NMCARS list styling needs.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Create the new "List 1" paragraph style (based on Heading 1) plus
#    the linked character styles Word auto-generates alongside it.
# ---------------------------------------------------------------------

$heading1 = $d.Styles("Heading1")

$list1 = $d.Styles.Add("List 1", 1)
$list1.BaseStyle = $d.Styles("Heading1")
$list1.ParagraphFormat.KeepWithNext = $false
$list1.ParagraphFormat.TabStops.ClearAll()
$list1.ParagraphFormat.TabStops.Add(184.3)
$list1.ParagraphFormat.SpaceBefore = 0
$list1.ParagraphFormat.SpaceAfter = 0
$list1.ParagraphFormat.Alignment = 0
$list1.ParagraphFormat.OutlineLevel = 10
$list1.Font.Bold = $false
$list1.Font.Size = 12

$heading1Char1 = $d.Styles.Add("Heading 1 Char1", 2)
$heading1Char1.BaseStyle = $d.Styles("DefaultParagraphFont")
$heading1Char1.Font.Bold = $true
$heading1Char1.Font.Size = 16

$list1Char = $d.Styles.Add("List 1 Char", 2)
$list1Char.BaseStyle = $d.Styles("Heading1Char1")
$list1Char.Font.Bold = $false
$list1Char.Font.Size = 12

# link Heading 1 <-> Heading 1 Char1, and List 1 <-> List 1 Char
$d.Styles("Heading1").LinkStyle = $d.Styles("Heading1Char1")
$d.Styles("Heading1Char1").LinkStyle = $d.Styles("Heading1")
$d.Styles("List1").LinkStyle = $d.Styles("List1Char")
$d.Styles("List1Char").LinkStyle = $d.Styles("List1")

# ---------------------------------------------------------------------
# 2. Add the "List 3" / "List 4" paragraph styles.
# ---------------------------------------------------------------------

$list3 = $d.Styles.Add("List 3", 1)
$list3.BaseStyle = $d.Styles("Normal")
$list3.ParagraphFormat.LeftIndent = 54
$list3.ParagraphFormat.FirstLineIndent = -18
$list3.NoSpaceBetweenParagraphsOfSameStyle = $true

$list4 = $d.Styles.Add("List 4", 1)
$list4.BaseStyle = $d.Styles("Normal")
$list4.ParagraphFormat.LeftIndent = 72
$list4.ParagraphFormat.FirstLineIndent = -18
$list4.NoSpaceBetweenParagraphsOfSameStyle = $true
$list4.UnhideWhenUsed = $true

# ---------------------------------------------------------------------
# 3. Apply the "List 1" style to the relevant paragraphs.
# ---------------------------------------------------------------------

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Unusual contract financing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.Style = "List 1"

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Submit requests for advance payments", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.Style = "List 1"

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("(c)(2) Information should be submitted", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.Style = "List 1"

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Contracting officers shall obtain", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.Style = "List 1"

# ---------------------------------------------------------------------
# 4. Tidy the built-in "List 2" style: drop its direct Courier New
#    run formatting.
# ---------------------------------------------------------------------

$list2 = $d.Styles("List2")
$list2.Font.Name = ""
